$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1573.5714
$ws.Range("I53").Value = 1305.2222
$ws.Range("J53").Value = 2056.6
$ws.Range("K53").Value = 1305.2222
$ws.Range("L53").Value = 2056.6
$ws.Range("M53").Value = -668.2221999999999
$ws.Range("N53").Value = -3330.6
$ws.Range("H70").Value = 2302.4443
$ws.Range("J70").Value = 2253.3333
$ws.Range("L70").Value = 6759.999899999999
$ws.Range("N70").Value = -7299.999899999999
$ws.Range("H73").Value = 2302.4443
$ws.Range("J73").Value = 2253.3333
$ws.Range("L73").Value = 6759.999899999999
$ws.Range("N73").Value = -8631.999899999999
$ws.Range("H107").Value = 423.44
$ws.Range("I107").Value = 399.875
$ws.Range("K107").Value = 399.875
$ws.Range("M107").Value = 1520.125
$ws.Range("H116").Value = 8925.333000000001
$ws.Range("I116").Value = 9000
$ws.Range("K116").Value = 9000
$ws.Range("M116").Value = -5558
$ws.Range("H132").Value = 2385.9
$ws.Range("I132").Value = 2127.3572
$ws.Range("K132").Value = 6382.071599999999
$ws.Range("M132").Value = -3852.071599999999
$ws.Range("H135").Value = 752.86957
$ws.Range("I135").Value = 681.5789
$ws.Range("K135").Value = 6134.2101
$ws.Range("M135").Value = -3599.2101
$ws.Range("H137").Value = 5083.923
$ws.Range("I137").Value = 5411.75
$ws.Range("K137").Value = 16235.25
$ws.Range("M137").Value = -13685.25
$ws.Range("H141").Value = 1800.1428
$ws.Range("I141").Value = 1734.9
$ws.Range("K141").Value = 5204.700000000001
$ws.Range("M141").Value = -24.70000000000073
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5568.4653
$ws.Range("I32").Value = 4685.7407
$ws.Range("K32").Value = 4685.7407
$ws.Range("M32").Value = -4398.7407
$ws.Range("H74").Value = 71510500
$ws.Range("I74").Value = 143018500
$ws.Range("K74").Value = 143018500
$ws.Range("M74").Value = -143017626
$ws.Range("H77").Value = 71510500
$ws.Range("I77").Value = 143018500
$ws.Range("K77").Value = 715092500
$ws.Range("M77").Value = -715088132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H94").Value = 1585.9166
$ws.Range("I94").Value = 547
$ws.Range("K94").Value = 547
$ws.Range("M94").Value = -96
$ws.Range("H134").Value = 3368.158
$ws.Range("I134").Value = 3474.5833
$ws.Range("K134").Value = 10423.7499
$ws.Range("M134").Value = -7888.749899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7902.1
$ws.Range("I31").Value = 4101.6
$ws.Range("K31").Value = 4101.6
$ws.Range("M31").Value = -3806.6
$ws.Range("H34").Value = 7902.1
$ws.Range("I34").Value = 4101.6
$ws.Range("K34").Value = 4101.6
$ws.Range("M34").Value = -3899.6
$ws.Range("H41").Value = 5500
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H62").Value = 3724.75
$ws.Range("I62").Value = 3724.75
$ws.Range("K62").Value = 3724.75
$ws.Range("M62").Value = -3100.75
$ws.Range("H65").Value = 3724.75
$ws.Range("I65").Value = 3724.75
$ws.Range("K65").Value = 18623.75
$ws.Range("M65").Value = -15503.75
$ws.Range("H107").Value = 3661.25
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 4631.6665
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 4631.6665
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -8471.666499999999
$ws.Range("H131").Value = 50697
$ws.Range("J131").Value = 50697
$ws.Range("L131").Value = 50697
$ws.Range("N131").Value = -60777
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5999.6665
$ws.Range("I5").Value = 2999.5
$ws.Range("J5").Value = 12000
$ws.Range("K5").Value = 8998.5
$ws.Range("L5").Value = 36000
$ws.Range("M5").Value = -8886.5
$ws.Range("N5").Value = -36224
$ws.Range("H12").Value = 554
$ws.Range("I12").Value = 251
$ws.Range("K12").Value = 753
$ws.Range("M12").Value = -580
$ws.Range("H34").Value = 182.61905
$ws.Range("I34").Value = 134.66667
$ws.Range("J34").Value = 302.5
$ws.Range("K34").Value = 404.00001
$ws.Range("L34").Value = 907.5
$ws.Range("M34").Value = -320.00001
$ws.Range("N34").Value = -1075.5
$ws.Range("H55").Value = 741.2727
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 942.3333
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 2826.9999
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -3180.9999
$ws.Range("H131").Value = 26945.818
$ws.Range("I131").Value = 95132.63
$ws.Range("J131").Value = 4216.879
$ws.Range("K131").Value = 285397.89
$ws.Range("L131").Value = 12650.637
$ws.Range("M131").Value = -280357.89
$ws.Range("N131").Value = -22730.637
$ws.Range("H132").Value = 2542.7778
$ws.Range("J132").Value = 3436.25
$ws.Range("L132").Value = 30926.25
$ws.Range("N132").Value = -35986.25
$ws.Range("H135").Value = 5999.6665
$ws.Range("I135").Value = 2999.5
$ws.Range("J135").Value = 12000
$ws.Range("K135").Value = 26995.5
$ws.Range("L135").Value = 108000
$ws.Range("M135").Value = -24460.5
$ws.Range("N135").Value = -113070
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H80").Value = 3718.2
$ws.Range("I80").Value = 3721.6
$ws.Range("K80").Value = 3721.6
$ws.Range("M80").Value = -2723.6
$ws.Range("H83").Value = 3718.2
$ws.Range("I83").Value = 3721.6
$ws.Range("K83").Value = 18608
$ws.Range("M83").Value = -13616
$ws.Range("H95").Value = 296137.66
$ws.Range("I95").Value = 200000
$ws.Range("J95").Value = 344206.5
$ws.Range("K95").Value = 200000
$ws.Range("L95").Value = 344206.5
$ws.Range("N95").Value = -349698.5
$ws.Range("M95").Value = -197254
$ws.Range("H102").Value = 2685.3635
$ws.Range("I102").Value = 1485.2858
$ws.Range("J102").Value = 4785.5
$ws.Range("K102").Value = 1485.2858
$ws.Range("L102").Value = 4785.5
$ws.Range("M102").Value = 136.7141999999999
$ws.Range("N102").Value = -8029.5
$ws.Range("H122").Value = 1896.3572
$ws.Range("I122").Value = 1659
$ws.Range("K122").Value = 4977
$ws.Range("M122").Value = -2527
$ws.Range("H123").Value = 39399.8
$ws.Range("J123").Value = 54999.5
$ws.Range("L123").Value = 54999.5
$ws.Range("N123").Value = -59899.5
$ws.Range("H132").Value = 4049.7273
$ws.Range("I132").Value = 4817
$ws.Range("J132").Value = 2003.6666
$ws.Range("K132").Value = 14451
$ws.Range("L132").Value = 6010.9998
$ws.Range("M132").Value = -11921
$ws.Range("N132").Value = -11070.9998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4089.8
$ws.Range("I7").Value = 3919.4
$ws.Range("J7").Value = 4146.6
$ws.Range("K7").Value = 3919.4
$ws.Range("L7").Value = 4146.6
$ws.Range("M7").Value = -3807.4
$ws.Range("N7").Value = -4370.6
$ws.Range("H9").Value = 171.8
$ws.Range("I9").Value = 171.8
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 171.8
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 52.19999999999999
$ws.Range("N9").ClearContents()
$ws.Range("H22").Value = 1074
$ws.Range("I22").Value = 536
$ws.Range("J22").Value = 2150
$ws.Range("K22").Value = 536
$ws.Range("L22").Value = 2150
$ws.Range("M22").Value = -241
$ws.Range("N22").Value = -2740
$ws.Range("H27").Value = 1074
$ws.Range("I27").Value = 536
$ws.Range("J27").Value = 2150
$ws.Range("K27").Value = 536
$ws.Range("L27").Value = 2150
$ws.Range("M27").Value = -429
$ws.Range("N27").Value = -2364
$ws.Range("H40").Value = 2959.4211
$ws.Range("I40").Value = 2987.7222
$ws.Range("K40").Value = 2987.7222
$ws.Range("M40").Value = -2851.7222
$ws.Range("H55").Value = 606.35
$ws.Range("I55").Value = 409
$ws.Range("K55").Value = 409
$ws.Range("M55").Value = -236
$ws.Range("H122").Value = 4809.5527
$ws.Range("I122").Value = 3929.1538
$ws.Range("K122").Value = 11787.4614
$ws.Range("M122").Value = -9337.4614
$ws.Range("H126").Value = 4089.8
$ws.Range("I126").Value = 3919.4
$ws.Range("J126").Value = 4146.6
$ws.Range("K126").Value = 11758.2
$ws.Range("L126").Value = 12439.8
$ws.Range("M126").Value = -9288.200000000001
$ws.Range("N126").Value = -17379.8
$ws.Range("H131").Value = 51743.668
$ws.Range("J131").Value = 51743.668
$ws.Range("L131").Value = 51743.668
$ws.Range("N131").Value = -61823.668
$ws.Range("H132").Value = 10207.986
$ws.Range("I132").Value = 7819.116
$ws.Range("K132").Value = 23457.348
$ws.Range("M132").Value = -20927.348
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 333346270
$ws.Range("I41").Value = 1000000000
$ws.Range("J41").Value = 19400
$ws.Range("K41").Value = 1000000000
$ws.Range("L41").Value = 19400
$ws.Range("M41").Value = -999999610
$ws.Range("N41").Value = -20180
$ws.Range("H126").Value = 5728
$ws.Range("I126").Value = 6487.636
$ws.Range("K126").Value = 19462.908
$ws.Range("M126").Value = -16992.908
$ws.Range("H132").Value = 3051.2222
$ws.Range("I132").Value = 3462.5334
$ws.Range("K132").Value = 10387.6002
$ws.Range("M132").Value = -7857.600199999999
